$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.675.66'
$ws.Range("E2").Value = '  +3.21%  '
$ws.Range("D3").Value = '2.198.43'
$ws.Range("E3").Value = '  +0.86%  '
$ws.Range("D5").Value = '''260.20'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.21%  '
$ws.Range("D6").Value = '''82.28'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +12.80%  '
$ws.Range("E7").Value = '  +2.90%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '''0.594'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.39%  '
$ws.Range("D10").Value = '''43.62'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +8.96%  '
$ws.Range("E11").Value = '  +1.26%  '
$ws.Range("D12").Value = '''6.97'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.28%  '
$ws.Range("E13").Value = '  +2.01%  '
$ws.Range("D14").Value = '2.525.51'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("E15").Value = '  +0.89%  '
$ws.Range("D16").Value = '2.187.12'
$ws.Range("E16").Value = '  +0.91%  '
$ws.Range("D17").Value = '''0.780'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.08%  '
$ws.Range("D18").Value = '43.595.42'
$ws.Range("E18").Value = '  +3.21%  '
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("D20").Value = '''69.81'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.05%  '
$ws.Range("E21").Value = '  +1.50%  '
$ws.Range("E22").Value = '  +16.09%  '
$ws.Range("D23").Value = '''230.69'
$ws.Range("D23").ClearFormats()
$ws.Range("D24").Value = '''8.86'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -4.78%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("E26").Value = '  +2.74%  '
$ws.Range("D27").Value = '''42.45'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +15.78%  '
$ws.Range("D28").Value = '''10.73'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +2.66%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '''2.24'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.45%  '
$ws.Range("E31").Value = '  +0.92%  '
$ws.Range("D32").Value = '''174.06'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.39%  '
$ws.Range("D33").Value = '''20.43'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.29%  '
$ws.Range("D34").Value = '''0.0874'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +6.90%  '
$ws.Range("E35").Value = '  +4.27%  '
$ws.Range("E36").Value = '  +6.69%  '
$ws.Range("E37").Value = '  +2.05%  '
$ws.Range("E38").Value = '  +6.74%  '
$ws.Range("D39").Value = '''0.0352'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +5.28%  '
$ws.Range("D40").Value = '''13.06'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +11.43%  '
$ws.Range("D41").Value = '''2.86'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +16.87%  '
$ws.Range("E42").Value = '  +2.88%  '
$ws.Range("D43").Value = '''64.39'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +8.79%  '
$ws.Range("D44").Value = '''5.47'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +6.47%  '
$ws.Range("E45").Value = '  +2.75%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '''0.0980'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''99.96'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.62%  '
$ws.Range("D48").Value = '''8.27'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.66%  '
$ws.Range("E49").Value = '  +4.13%  '
$ws.Range("E50").Value = '  +2.86%  '
$ws.Range("D51").Value = '''0.440'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.80%  '
